# Auto-generated Excel COM-interop script applying market-data refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (Table_* ranges).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1855.6666
$ws.Range("I12").Value = 1855.6666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1855.6666
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1685.6666
$ws.Range("N12").ClearContents()

$ws.Range("H42").Value = 308.85715
$ws.Range("I42").Value = 30.6
$ws.Range("K42").Value = 91.80000000000001
$ws.Range("M42").Value = 138.2

$ws.Range("H43").Value = 4302.355
$ws.Range("I43").Value = 4074.35
$ws.Range("J43").Value = 4716.909
$ws.Range("K43").Value = 4074.35
$ws.Range("L43").Value = 4716.909
$ws.Range("M43").Value = -4005.35
$ws.Range("N43").Value = -4854.909

$ws.Range("H86").Value = 8161.143
$ws.Range("J86").Value = 8082
$ws.Range("L86").Value = 8082
$ws.Range("N86").Value = -10328

$ws.Range("H89").Value = 8161.143
$ws.Range("J89").Value = 8082
$ws.Range("L89").Value = 40410
$ws.Range("N89").Value = -51642

$ws.Range("H98").Value = 2175
$ws.Range("I98").Value = 2175
$ws.Range("K98").Value = 2175
$ws.Range("M98").Value = -677

$ws.Range("H107").Value = 530
$ws.Range("I107").Value = 570.0714
$ws.Range("K107").Value = 570.0714
$ws.Range("M107").Value = 1349.9286

$ws.Range("H122").Value = 2175
$ws.Range("I122").Value = 2175
$ws.Range("K122").Value = 6525
$ws.Range("M122").Value = -4075

$ws.Range("H125").Value = 4232.8335
$ws.Range("I125").Value = 2215.5
$ws.Range("J125").Value = 8267.5
$ws.Range("K125").Value = 19939.5
$ws.Range("L125").Value = 74407.5
$ws.Range("M125").Value = -17479.5
$ws.Range("N125").Value = -79327.5

$ws.Range("H137").Value = 1511.9445
$ws.Range("I137").Value = 1425.7273
$ws.Range("K137").Value = 4277.1819
$ws.Range("M137").Value = -1727.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 13082.714
$ws.Range("I28").Value = 13082.714
$ws.Range("K28").Value = 13082.714
$ws.Range("M28").Value = -12890.714

$ws.Range("H99").Value = 13082.714
$ws.Range("I99").Value = 13082.714
$ws.Range("K99").Value = 13082.714
$ws.Range("M99").Value = -10087.714

$ws.Range("H122").Value = 2432.2
$ws.Range("I122").Value = 1219
$ws.Range("J122").Value = 3645.4
$ws.Range("K122").Value = 3657
$ws.Range("L122").Value = 10936.2
$ws.Range("M122").Value = -1207
$ws.Range("N122").Value = -15836.2

$ws.Range("H132").Value = 1679.0667
$ws.Range("I132").Value = 1602.7693
$ws.Range("K132").Value = 4808.3079
$ws.Range("M132").Value = -2278.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2596.9285
$ws.Range("I86").Value = 1816.5
$ws.Range("J86").Value = 3637.5
$ws.Range("K86").Value = 1816.5
$ws.Range("L86").Value = 3637.5
$ws.Range("M86").Value = -693.5
$ws.Range("N86").Value = -5883.5

$ws.Range("H89").Value = 2596.9285
$ws.Range("I89").Value = 1816.5
$ws.Range("J89").Value = 3637.5
$ws.Range("K89").Value = 9082.5
$ws.Range("L89").Value = 18187.5
$ws.Range("M89").Value = -3466.5
$ws.Range("N89").Value = -29419.5

$ws.Range("H107").Value = 925
$ws.Range("I107").Value = 999.6667
$ws.Range("K107").Value = 999.6667
$ws.Range("M107").Value = 920.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11483553
$ws.Range("I86").Value = 11483553
$ws.Range("K86").Value = 11483553
$ws.Range("M86").Value = -11482430

$ws.Range("H89").Value = 11483553
$ws.Range("I89").Value = 11483553
$ws.Range("K89").Value = 57417765
$ws.Range("M89").Value = -57412149

$ws.Range("H132").Value = 2328.742
$ws.Range("I132").Value = 2392.2964
$ws.Range("K132").Value = 7176.889200000001
$ws.Range("M132").Value = -4646.889200000001

$ws.Range("H134").Value = 2067.6365
$ws.Range("I134").Value = 1774.4
$ws.Range("K134").Value = 5323.200000000001
$ws.Range("M134").Value = -2788.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 239.25
$ws.Range("I47").Value = 276.5
$ws.Range("K47").Value = 829.5
$ws.Range("M47").Value = -398.5

$ws.Range("H113").Value = 763.36365
$ws.Range("J113").Value = 785.4286
$ws.Range("L113").Value = 2356.2858
$ws.Range("N113").Value = -6696.2858

$ws.Range("H139").Value = 2563.6365
$ws.Range("I139").Value = 2950
$ws.Range("J139").Value = 2100
$ws.Range("K139").Value = 8850
$ws.Range("L139").Value = 6300
$ws.Range("M139").Value = -3710
$ws.Range("N139").Value = -16580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 221.90909
$ws.Range("I2").Value = 171.25
$ws.Range("J2").Value = 250.85715
$ws.Range("K2").Value = 171.25
$ws.Range("L2").Value = 250.85715
$ws.Range("M2").Value = -58.25
$ws.Range("N2").Value = -476.85715

$ws.Range("H80").Value = 5430.143
$ws.Range("I80").Value = 3331.6667
$ws.Range("J80").Value = 7004
$ws.Range("K80").Value = 3331.6667
$ws.Range("L80").Value = 7004
$ws.Range("M80").Value = -2333.6667
$ws.Range("N80").Value = -9000

$ws.Range("H83").Value = 5430.143
$ws.Range("I83").Value = 3331.6667
$ws.Range("J83").Value = 7004
$ws.Range("K83").Value = 16658.3335
$ws.Range("L83").Value = 35020
$ws.Range("M83").Value = -11666.3335
$ws.Range("N83").Value = -45004

$ws.Range("H102").Value = 883.25
$ws.Range("I102").Value = 687.7273
$ws.Range("K102").Value = 687.7273
$ws.Range("M102").Value = 934.2727

$ws.Range("H122").Value = 4089.1
$ws.Range("I122").Value = 4235.375
$ws.Range("J122").Value = 3504
$ws.Range("K122").Value = 12706.125
$ws.Range("L122").Value = 10512
$ws.Range("M122").Value = -10256.125
$ws.Range("N122").Value = -15412

$ws.Range("H126").Value = 4097.6665
$ws.Range("I126").Value = 3998.5
$ws.Range("K126").Value = 11995.5
$ws.Range("M126").Value = -9525.5

$ws.Range("H132").Value = 2525.9
$ws.Range("I132").Value = 2728.889
$ws.Range("J132").Value = 699
$ws.Range("K132").Value = 8186.667
$ws.Range("L132").Value = 2097
$ws.Range("M132").Value = -5656.667
$ws.Range("N132").Value = -7157

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2714.8572
$ws.Range("I40").Value = 2714.8572
$ws.Range("K40").Value = 2714.8572
$ws.Range("M40").Value = -2578.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 75644
$ws.Range("J110").Value = 75644
$ws.Range("L110").Value = 75644
$ws.Range("N110").Value = -83824

$ws.Range("H126").Value = 1637.4
$ws.Range("I126").Value = 1637.4
$ws.Range("K126").Value = 4912.200000000001
$ws.Range("M126").Value = -2442.200000000001

$ws.Range("H132").Value = 1756.8462
$ws.Range("I132").Value = 1756.8462
$ws.Range("K132").Value = 5270.5386
$ws.Range("M132").Value = -2740.5386

$ws.Range("H136").Value = 703
$ws.Range("I136").Value = 671.12
$ws.Range("K136").Value = 2013.36
$ws.Range("M136").Value = 536.6399999999999

Write-Host "Applied 179 cell updates and 1 clears across 8 sheets."